$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Truth-table input combinations (A,B,C,D,E) for rows 75-82, continuing the
# existing A=0/B=1 block (rows 67-74) with the A=1/B=1 block.
$inputs = @(
    @(1,1,0,0,0),
    @(1,1,0,1,0),
    @(1,1,0,0,1),
    @(1,1,0,1,1),
    @(1,1,1,0,0),
    @(1,1,1,1,0),
    @(1,1,1,0,1),
    @(1,1,1,1,1)
)

$row = 75
foreach ($vals in $inputs) {
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]

    $ws.Cells.Item($row, 6).Formula = "=B$row*2 + E$row"
    $ws.Cells.Item($row, 7).Formula = "=A$row * 2 + D$row"
    $ws.Cells.Item($row, 8).Formula = "=D$row"
    $ws.Cells.Item($row, 9).Formula = "=E$row"
    $ws.Cells.Item($row, 10).Formula = "=F$row"
    $ws.Cells.Item($row, 11).Formula = "=C$row * 2 + D$row"
    $ws.Cells.Item($row, 12).Formula = "=D$row"
    $ws.Cells.Item($row, 13).Formula = "=E$row"
    $ws.Cells.Item($row, 14).Formula = "=C$row*2 + E$row"
    $ws.Cells.Item($row, 15).Formula = "=G$row"

    $row++
}

$excel.ActiveWindow.ScrollRow = 40
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("D75").Select()
